# Adds a new "11-ago" day column (column AX) to Sheet1, mirroring the
# existing "10-ago" column (AW): a text-formatted header in row 1 and
# numeric values in rows 2-18. Also updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column is the one right after the current last column (AW = 49 -> AX = 50)
$newCol = 50

# Header cell: same "text" number format as the other date headers (e.g. AW1)
$headerCell = $ws.Cells.Item(1, $newCol)
$headerCell.NumberFormat = "@"
$headerCell.Value2 = "11-ago"

# Data values for rows 2-18 (mirrors the new "11-ago" column in the diff)
$values = @(
    0,
    16.37161535676869,
    19.189888191091537,
    15.51038477889068,
    0,
    22.283401534305217,
    15.622993165001118,
    14.050738871089072,
    11.625658964831546,
    9.5946551562429665,
    0,
    13.042864421251565,
    0,
    0,
    8.1121270278387243,
    0,
    0
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, $newCol).Value2 = $values[$i]
}

# Update the selected cell shown in the saved view (AY7 -> AY3)
$ws.Range("AY3").Select()
